# Apply the error-log edit described by the commit:
# "add animation for chatbot, make new error log for windows error"
#
# - Rename all occurrences of user "Yoko Suzuki" -> "Taro Fujita" (C2:C16)
# - Renumber capture image paths in column J (bdot20240415_141953/... -> bdot20240415_141954/N.png)
# - Rewrite the explanation text (column K) to describe a fresh scenario where
#   rows 5-7 now hold the new Windows error (0x80240fff) instead of the old one (0x80244007)
# - Move the error_type / error_content (L/M) from row 7 to row 5, clearing them on row 7
# - Swap B5/B7 type values ("operation"/"error") to match the relocated error row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Cell($addr, $value) {
    $ws.Range($addr).Value = $value
}

# --- Column C: user_name, rows 2-16 ---
foreach ($r in 2..16) {
    Set-Cell "C$r" "Taro Fujita"
}

# --- Column J: capimg, rows 2-16 ---
Set-Cell "J2"  "bdot20240415_141954/1.png"
Set-Cell "J3"  "bdot20240415_141954/2.png"
Set-Cell "J4"  "bdot20240415_141954/3.png"
Set-Cell "J5"  "bdot20240415_141954/4.png"
Set-Cell "J6"  "bdot20240415_141954/5.png"
Set-Cell "J7"  "bdot20240415_141954/5.png"
Set-Cell "J8"  "bdot20240415_141954/6.png"
Set-Cell "J9"  "bdot20240415_141954/7.png"
Set-Cell "J10" "bdot20240415_141954/8.png"
Set-Cell "J11" "bdot20240415_141954/9.png"
Set-Cell "J12" "bdot20240415_141954/10.png"
Set-Cell "J13" "bdot20240415_141954/1.png"
Set-Cell "J14" "bdot20240415_141954/2.png"
Set-Cell "J15" "bdot20240415_141954/3.png"
Set-Cell "J16" "bdot20240415_141954/11.png"

# --- Column B: type, rows 5 & 7 swap ---
Set-Cell "B5" "error"
Set-Cell "B7" "operation"

# --- Column K: explanation, rows 2-16 ---
Set-Cell "K2"  "「スタート」ボタンをクリックする"
Set-Cell "K3"  "メニューから「設定」アイコンをクリックする"
Set-Cell "K4"  "左側のメニューからWindows Updateをクリックし、Windows Update画面に移動する"
Set-Cell "K5"  "0x80240fff エラー"
Set-Cell "K6"  "デスクトップ画面の左下にある「スタート」ボタンを右クリックする"
Set-Cell "K7"  "メニューからターミナル(管理者)をクリックする"
Set-Cell "K8"  "ユーザーアカウント制御と表示されているウィンドウが開いたことを確認する"
Set-Cell "K9"  "PowerShellウィンドウに start-transcript と入力し、[Enter]キーを押す"
Set-Cell "K10" "wuauclt.exe /resetauthorization /detectnow と入力し、[Enter]キーを押す"
Set-Cell "K11" "netsh winhttp show proxy と入力し、[Enter]キーを押す"
Set-Cell "K12" "netsh winhttp reset proxy と入力し、[Enter]キーを押す"
Set-Cell "K13" "「スタート」ボタンをクリックする"
Set-Cell "K14" "メニューから「設定」アイコンをクリックする"
Set-Cell "K15" "左側のメニューからWindows Updateをクリックし、Windows Update画面に移動する"
Set-Cell "K16" "「更新プログラムのチェック」ボタンをクリックする"

# --- Column L/M: error_type / error_content move from row 7 to row 5 ---
Set-Cell "L5" "Error W"
Set-Cell "M5" " エラーの Windows"
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
